$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Function Index"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Function Index")

# Reorder / update "Calls" lists (text only changes, same members)
$ws1.Range("D2").Value = "init, evaluate, nowMs, setAspect"

$ws1.Range("D3").Value = "setAspect, nowMs, isOccupied, update, evaluateControllerLogic, isHealthy"
$ws1.Range("F3").Value = 3

$ws1.Range("D6").Value = "writeLamp, simultaneously"

$ws1.Range("D7").Value = "configure, readRawClear"

$ws1.Range("D10").Value = "digitalWrite, else"

# Row 14 now describes the new helper function (was "evaluate")
$ws1.Range("A14").Value = "computeControllerFresh"
$ws1.Range("B14").Value = "src\logic\ControllerHelpers.cpp"

# New row 15: evaluateControllerLogic
$ws1.Range("A15").Value = "evaluateControllerLogic"
$ws1.Range("B15").Value = "src\logic\ControllerLogic.cpp"
$ws1.Range("C15").Value = "cpp"
$ws1.Range("D15").Value = "evaluate, computeControllerFresh"
$ws1.Range("E15").Value = $false
$ws1.Range("F15").Value = 2

# New row 16: evaluate (moved back to Interlocking.cpp, no calls)
$ws1.Range("A16").Value = "evaluate"
$ws1.Range("B16").Value = "src\logic\Interlocking.cpp"
$ws1.Range("C16").Value = "cpp"
$ws1.Range("E16").Value = $false
$ws1.Range("F16").Value = 1

# ---------------------------------------------------------------------------
# Sheet "Call Graph"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Call Graph")

$ws2.Range("B2").Value = "init, evaluate, nowMs, setAspect"
$ws2.Range("B3").Value = "setAspect, nowMs, isOccupied, update, evaluateControllerLogic, isHealthy"
$ws2.Range("B6").Value = "writeLamp, simultaneously"
$ws2.Range("B7").Value = "configure, readRawClear"
$ws2.Range("B10").Value = "digitalWrite, else"

$ws2.Range("A14").Value = "computeControllerFresh"

$ws2.Range("A15").Value = "evaluateControllerLogic"
$ws2.Range("B15").Value = "evaluate, computeControllerFresh"

$ws2.Range("A16").Value = "evaluate"

# ---------------------------------------------------------------------------
# Sheet "File Summaries"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("File Summaries")

$ws3.Range("E2").Value = 3

# Insert two blank rows before the old row 9 (src\logic\Interlocking.cpp),
# pushing it down to row 11, then populate rows 9 and 10 with the new files.
$ws3.Range("A9:A10").EntireRow.Insert()

$ws3.Range("A9").Value = "src\logic\ControllerHelpers.cpp"
$ws3.Range("B9").Value = 1
$ws3.Range("C9").Value = $true
$ws3.Range("D9").Value = $true
$ws3.Range("E9").Value = 1

$ws3.Range("A10").Value = "src\logic\ControllerLogic.cpp"
$ws3.Range("B10").Value = 1
$ws3.Range("C10").Value = $true
$ws3.Range("D10").Value = $true
$ws3.Range("E10").Value = 2
$ws3.Range("F10").Value = "src\logic\Interlocking.cpp, src\logic\ControllerHelpers.cpp"

# ---------------------------------------------------------------------------
# Sheet "Class Roles" (reordered rows)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Class Roles")

$ws4.Range("A2").Value = "TrackCircuitInput"
$ws4.Range("B2").Value = "MIXED"

$ws4.Range("A3").Value = "BlockController"
$ws4.Range("B3").Value = "MIXED"

$ws4.Range("A4").Value = "ArduinoGpio"
$ws4.Range("B4").Value = "HARDWARE"

$ws4.Range("A5").Value = "MockGpio"
$ws4.Range("B5").Value = "HARDWARE"

$ws4.Range("A6").Value = "SignalHead"
$ws4.Range("B6").Value = "MIXED"
